# Adapt column header formatting to respective input file names:
#   <name>_old -> <name>_FV2404   (columns A..J)
#   <name>_new -> <name>_FV2410   (columns L..U)
# then wrap the used range in an Excel Table (Table1) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$colsOld = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$colsNew = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($colsOld[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($colsNew[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# --- Create the Excel Table over A1:U78 --------------------------------
# The header row already carries manual formatting (bold/fill/border) from
# the original sheet. If that formatting is present when the table is
# created, Excel bakes it into a header-row dxf/style override. The source
# workbook doesn't have that override, so stash the header formatting,
# clear it, create the table, then restore the formatting byte-for-byte.

$headerRange = $ws.Range("A1:U1")
$scratchRow = 1000
$scratch = $ws.Range("A" + $scratchRow + ":U" + $scratchRow)

[void]$headerRange.Copy()
[void]$scratch.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

[void]$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

[void]$scratch.Copy()
[void]$headerRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
[void]$ws.Rows.Item($scratchRow).Delete()

# --- Freeze the header row ----------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
